# Apply the table style change to the three tables (slides 14, 15, 16)
# that currently use the custom "Table_0" style and should instead use
# the built-in table style {8328220B-B384-473B-954A-A5236DE1C7D1}.
$p = $ppt.ActivePresentation

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle("{8328220B-B384-473B-954A-A5236DE1C7D1}")
        }
    }
}

# Re-colour the presentation's theme (theme1.xml, used by the slide
# master/design "Integral") so its 12 scheme colours become the
# standard default "Office" palette - matching the theme that the
# notes master already uses. (RGB values are passed in BGR/VBA long
# colour order, as required by ThemeColorScheme.Item(n).RGB.)
$slide1 = $p.Slides.Item(1)
$colours = $slide1.ThemeColorScheme

$colours.Item(1).RGB  = 0x000000  # dk1
$colours.Item(2).RGB  = 0xFFFFFF  # lt1
$colours.Item(3).RGB  = 0x6A5444  # dk2      -> 44546A
$colours.Item(4).RGB  = 0xE6E6E7  # lt2      -> E7E6E6
$colours.Item(5).RGB  = 0xD59B5B  # accent1  -> 5B9BD5
$colours.Item(6).RGB  = 0x317DED  # accent2  -> ED7D31
$colours.Item(7).RGB  = 0xA5A5A5  # accent3  -> A5A5A5
$colours.Item(8).RGB  = 0x00C0FF  # accent4  -> FFC000
$colours.Item(9).RGB  = 0xC47244  # accent5  -> 4472C4
$colours.Item(10).RGB = 0x47AD70  # accent6  -> 70AD47
$colours.Item(11).RGB = 0xC16305  # hlink    -> 0563C1
$colours.Item(12).RGB = 0x724F95  # folHlink -> 954F72
